$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row of results (AdaBoost run) below the existing data.
$ws.Range("B26").Value = "AdaBoost"
$ws.Range("C26").Value = "c2"
$ws.Range("D26").Value = "k5"
$ws.Range("E26").Value = "est 300"
$ws.Range("F26").Value = "depth 3"
$ws.Range("H26").Value = 0.9299

# Update the active selection to match the post-edit state.
$ws.Range("G26").Select()
